# Update "想去人数" (want-to-go count) figures in column F across the four
# sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 60
$ws.Range("F8").Value = 311
$ws.Range("F11").Value = 10537
$ws.Range("F15").Value = 2005
$ws.Range("F16").Value = 883
$ws.Range("F17").Value = 29
$ws.Range("F18").Value = 6
$ws.Range("F20").Value = 58
$ws.Range("F23").Value = 110
$ws.Range("F24").Value = 171
$ws.Range("F25").Value = 667
$ws.Range("F27").Value = 204
$ws.Range("F28").Value = 2360
$ws.Range("F29").Value = 654
$ws.Range("F30").Value = 3017
$ws.Range("F31").Value = 989
$ws.Range("F36").Value = 909
$ws.Range("F37").Value = 14
$ws.Range("F38").Value = 17
$ws.Range("F41").Value = 1168
$ws.Range("F44").Value = 122
$ws.Range("F45").Value = 214
$ws.Range("F46").Value = 46
$ws.Range("F47").Value = 9

# --- 演出 (Performances) ---------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 6
$ws.Range("F14").Value = 248

# --- 本地生活 (Local life) -------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 402

# --- 全部类型 (All types) --------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 402
$ws.Range("F5").Value = 6
$ws.Range("F8").Value = 60
$ws.Range("F11").Value = 311
$ws.Range("F13").Value = 10537
$ws.Range("F17").Value = 2005
$ws.Range("F18").Value = 883
$ws.Range("F19").Value = 29
$ws.Range("F20").Value = 6
$ws.Range("F24").Value = 110
$ws.Range("F25").Value = 171
$ws.Range("F28").Value = 667
$ws.Range("F30").Value = 204
$ws.Range("F31").Value = 2360
$ws.Range("F32").Value = 654
$ws.Range("F33").Value = 3017
$ws.Range("F34").Value = 989
$ws.Range("F36").Value = 909
$ws.Range("F37").Value = 17
$ws.Range("F40").Value = 1168
$ws.Range("F43").Value = 122
$ws.Range("F44").Value = 214
$ws.Range("F46").Value = 9
